# Daily attendance processing - 2026-01-28 14:03:37
# Swap the order of the "Recorded By" entries (column G) so that the
# dnasr281@gmail.com entry appears first, e.g.
#   "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
#   "admin@admin.com, dnasr281@gmail.com" -> "dnasr281@gmail.com, admin@admin.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($null -ne $val -and $val -is [string] -and $val.Contains(", ")) {
        $parts = $val.Split(", ")
        if ($parts.Count -eq 2 -and $parts[1] -eq "dnasr281@gmail.com") {
            $cell.Value2 = $parts[1] + ", " + $parts[0]
        }
    }
}
